# Update the expected-result strings in the "TestData" sheet (column E,
# rows 4-6 then 3) so that the shared-string table gains the three new
# messages in the same order they were newly introduced, replacing the old
# single placeholder string ("กรุณากรอกข้อมูลให้ถูกต้อง") that is no longer
# referenced anywhere.
$wb = $excel.ActiveWorkbook

$wsTestData = $wb.Worksheets.Item("TestData")

$wsTestData.Range("E4").Value = '"เลขรหัสประจำตัว 13 หลักไม่ถูกต้อง"'
$wsTestData.Range("E5").Value = '"เลขรหัสประจำตัว 13 หลักไม่ถูกต้อง"'
$wsTestData.Range("E6").Value = '"กรุณากรอกข้อมูล?"'
$wsTestData.Range("E3").Value = '"ไม่พบประเภทบัตรบุคคลไร้รัฐไร้สัญชาติ"'

# Restore the TestStep sheet's last-used selection before switching away
# from it, then make TestData the active sheet/tab with its own selection.
$wsTestStep = $wb.Worksheets.Item("TestStep")
$wsTestStep.Activate()
$wsTestStep.Range("F11").Select()

$wsTestData.Activate()
$wsTestData.Range("D7").Select()
